# Insert a new weekly record at the top of the data (row 4), pushing all
# existing data rows down by one. The sheet already has a constant header
# block in rows 1-3 (title/labels) and the first three data rows share the
# same Mercado/Region/Category metadata as every other row, so we only need
# to populate the columns that vary per record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4:58 down to 5:59, inserting a fresh blank row 4.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the latest market reading.
$ws.Cells.Item(4, 1).Value = 11
$ws.Cells.Item(4, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value = "Bíobío"
$ws.Cells.Item(4, 4).Value = 44630
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 100112012
$ws.Cells.Item(4, 7).Value = "Espinaca"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 150
$ws.Cells.Item(4, 11).Value = 10000
$ws.Cells.Item(4, 12).Value = 10000
$ws.Cells.Item(4, 13).Value = 10000
$ws.Cells.Item(4, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(4, 15).Value = "Región Metropolitana"
$ws.Cells.Item(4, 16).Value = 1000
$ws.Cells.Item(4, 17).Value = 10
$ws.Cells.Item(4, 18).Value = "Hortaliza"
